# soil weights... need to work on
# Add new "dry.wt reading 2" / percent-loss data to the soil.moisture sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("soil.moisture")

# --- New column F (second dry weight reading) for rows 2-11 ---
$ws.Range("F2").Value = 4.0374999999999996
$ws.Range("F3").Value = 4.4884000000000004
$ws.Range("F4").Value = 4.6595000000000004
$ws.Range("F5").Value = 4.4255000000000004
$ws.Range("F6").Value = 4.9896000000000003
$ws.Range("F7").Value = 5.0781999999999998
$ws.Range("F8").Value = 4.6997999999999998
$ws.Range("F9").Value = 4.5034999999999998
$ws.Range("F10").Value = 4.4480000000000004
$ws.Range("F11").Value = 4.7872000000000003

# --- New column I (second day/date label) for rows 2-11 ---
# Create the "20.FEB" shared string before the "percent" header string so
# the shared string table ends up with 20.FEB then percent, in that order.
$ws.Range("I2").Value = "20.FEB"
$ws.Range("I3").Value = "20.FEB"
$ws.Range("I4").Value = "20.FEB"
$ws.Range("I5").Value = "20.FEB"
$ws.Range("I6").Value = "20.FEB"
$ws.Range("I7").Value = "20.FEB"
$ws.Range("I8").Value = "20.FEB"
$ws.Range("I9").Value = "20.FEB"
$ws.Range("I10").Value = "20.FEB"
$ws.Range("I11").Value = "20.FEB"

# --- New column J (second time reading) for rows 2-11 ---
$ws.Range("J2").Value = 0.45833333333333331
$ws.Range("J3").Value = 0.45833333333333331
$ws.Range("J4").Value = 0.45833333333333331
$ws.Range("J5").Value = 0.45833333333333331
$ws.Range("J6").Value = 0.45833333333333331
$ws.Range("J7").Value = 0.45833333333333331
$ws.Range("J8").Value = 0.45833333333333331
$ws.Range("J9").Value = 0.45833333333333331
$ws.Range("J10").Value = 0.45833333333333331
$ws.Range("J11").Value = 0.45833333333333331

# Use the same time format Excel used for column H
$ws.Range("J2:J11").NumberFormat = $ws.Range("H2").NumberFormat

# --- New column K (total hrs reading 2) for rows 2-11 ---
$ws.Range("K2").Value = 95.5
$ws.Range("K3").Value = 95.5
$ws.Range("K4").Value = 95.5
$ws.Range("K5").Value = 95.5
$ws.Range("K6").Value = 95.5
$ws.Range("K7").Value = 95.5
$ws.Range("K8").Value = 95.5
$ws.Range("K9").Value = 95.5
$ws.Range("K10").Value = 95.5
$ws.Range("K11").Value = 95.5

# --- New column L: percent moisture lost, header + formula ---
$ws.Range("L1").Value = "percent"
$ws.Range("L2").Formula = "=((C2-F2)/C2)*100"
$ws.Range("L3:L11").Formula = "=((C3-F3)/C3)*100"

# Update the selection on the soil.moisture sheet and make it the active tab
$ws.Range("O12").Select()
$ws.Activate()
